$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.142437902037841
$ws.Range("D2").Value = 0.2436499207375107
$ws.Range("E2").Value = 0.1795842642411571
$ws.Range("F2").Value = 1.037516806086579
$ws.Range("G2").Value = 0.4977741487897163
$ws.Range("H2").Value = 0.6405293565296475
$ws.Range("J2").Value = 0.1779080579731058
$ws.Range("K2").Value = 0.6062153509819268
$ws.Range("N2").Value = 1.332010824900088
$ws.Range("O2").Value = 2.239954461524746
$ws.Range("B3").Value = 0.1329322490822022
$ws.Range("D3").Value = 0.2376561938734341
$ws.Range("E3").Value = 0.1752894354209715
$ws.Range("F3").Value = 1.035818188712852
$ws.Range("G3").Value = 0.4973056651683194
$ws.Range("H3").Value = 0.6439431398294246
$ws.Range("J3").Value = 0.1737421781387312
$ws.Range("K3").Value = 0.5302917244039236
$ws.Range("N3").Value = 1.34135563581652
$ws.Range("O3").Value = 2.245846444442208
$ws.Range("B4").Value = 0.1271636699005683
$ws.Range("D4").Value = 0.2340765145819859
$ws.Range("E4").Value = 0.1727412937730257
$ws.Range("F4").Value = 1.035359426712112
$ws.Range("G4").Value = 0.4973646131627447
$ws.Range("H4").Value = 0.6463244712575715
$ws.Range("J4").Value = 0.1712845774738625
$ws.Range("K4").Value = 0.4835337601548702
$ws.Range("N4").Value = 1.347579869274604
$ws.Range("O4").Value = 2.250786786596194
$ws.Range("B5").Value = 0.1248302004291304
$ws.Range("D5").Value = 0.2326431659262198
$ws.Range("E5").Value = 0.1717253228349449
$ws.Range("F5").Value = 1.035319460365407
$ws.Range("G5").Value = 0.4974757318285015
$ws.Range("H5").Value = 0.6473666664717541
$ws.Range("J5").Value = 0.1703083277189066
$ws.Range("K5").Value = 0.4644454815687027
$ws.Range("N5").Value = 1.350238794108506
$ws.Range("O5").Value = 2.253132587522842
$ws.Range("B6").Value = 0.1244437785083647
$ws.Range("D6").Value = 0.2324066973731078
$ws.Range("E6").Value = 0.1715579774382618
$ws.Range("F6").Value = 1.035321703809629
$ws.Range("G6").Value = 0.4974994417305254
$ws.Range("H6").Value = 0.6475440593136952
$ws.Range("J6").Value = 0.170147747647988
$ws.Range("K6").Value = 0.4612738709025734
$ws.Range("N6").Value = 1.350687710232982
$ws.Range("O6").Value = 2.25354219198698
$ws.Range("B7").Value = 0.1271321297210335
$ws.Range("D7").Value = 0.2340570809259788
$ws.Range("E7").Value = 0.1727275011827132
$ws.Range("F7").Value = 1.035358292475621
$ws.Range("G7").Value = 0.4973657591737961
$ws.Range("H7").Value = 0.6463382359477663
$ws.Range("J7").Value = 0.171271309182643
$ws.Range("K7").Value = 0.483276464992656
$ws.Range("N7").Value = 1.347615232243264
$ws.Range("O7").Value = 2.250817076343935
$ws.Range("B8").Value = 0.1391463659018939
$ws.Range("D8").Value = 0.2415625019762757
$ws.Range("E8").Value = 0.178084991875231
$ws.Range("F8").Value = 1.036809918806831
$ws.Range("G8").Value = 0.4975406423832638
$ws.Range("H8").Value = 0.6416472560873103
$ws.Range("J8").Value = 0.176450868447624
$ws.Range("K8").Value = 0.5800668211767004
$ws.Range("N8").Value = 1.335132047242418
$ws.Range("O8").Value = 2.241711487854261
$ws.Range("B9").Value = 0.163238089964068
$ws.Range("D9").Value = 0.2570732138505889
$ws.Range("E9").Value = 0.189294365655023
$ws.Range("F9").Value = 1.044289874208246
$ws.Range("G9").Value = 0.5006371510494745
$ws.Range("H9").Value = 0.6347095393731905
$ws.Range("J9").Value = 0.1874029710348992
$ws.Range("K9").Value = 0.7687082869792903
$ws.Range("N9").Value = 1.314505287635257
$ws.Range("O9").Value = 2.234353269250136
$ws.Range("B10").Value = 0.1812547692589561
$ws.Range("D10").Value = 0.2689471111230972
$ws.Range("E10").Value = 0.1979570000670918
$ws.Range("F10").Value = 1.05260965034968
$ws.Range("G10").Value = 0.5045967847211585
$ws.Range("H10").Value = 0.6309885342564883
$ws.Range("J10").Value = 0.1959344659256317
$ws.Range("K10").Value = 0.9065407090829467
$ws.Range("N10").Value = 1.301689961491959
$ws.Range("O10").Value = 2.235354996549574
$ws.Range("B11").Value = 0.189518208776974
$ws.Range("D11").Value = 0.2744517020786503
$ws.Range("E11").Value = 0.2019903044413596
$ws.Range("F11").Value = 1.057007869219277
$ws.Range("G11").Value = 0.5067652787718657
$ws.Range("H11").Value = 0.6295940947780281
$ws.Range("J11").Value = 0.1999211061769302
$ws.Range("K11").Value = 0.9690682674386153
$ws.Range("N11").Value = 1.296365913117739
$ws.Range("O11").Value = 2.237204068830522
$ws.Range("B12").Value = 0.1926569009091139
$ws.Range("D12").Value = 0.27655085845835
$ws.Range("E12").Value = 0.2035308824979509
$ws.Range("F12").Value = 1.058761540920344
$ws.Range("G12").Value = 0.507639324183927
$ws.Range("H12").Value = 0.629108903668822
$ws.Range("J12").Value = 0.2014459168061933
$ws.Range("K12").Value = 0.9927197624367068
$ws.Range("N12").Value = 1.294422405923292
$ws.Range("O12").Value = 2.238104727082202
$ws.Range("B13").Value = 0.1919805079010075
$ws.Range("D13").Value = 0.2760981165514664
$ws.Range("E13").Value = 0.2031985030742121
$ws.Range("F13").Value = 1.058379936303453
$ws.Range("G13").Value = 0.5074487296537882
$ws.Range("H13").Value = 0.6292114929536012
$ws.Range("J13").Value = 0.2011168479750722
$ws.Range("K13").Value = 0.987627183986632
$ws.Range("N13").Value = 1.294837748286369
$ws.Range("O13").Value = 2.237901837554261
$ws.Range("B14").Value = 0.1897762415301969
$ws.Range("D14").Value = 0.2746241073730715
$ws.Range("E14").Value = 0.2021167834942119
$ws.Range("F14").Value = 1.057150378304854
$ws.Range("G14").Value = 0.5068361267153705
$ws.Range("H14").Value = 0.629553319190947
$ws.Range("J14").Value = 0.2000462496895494
$ws.Range("K14").Value = 0.9710146262027308
$ws.Range("N14").Value = 1.296204565420965
$ws.Range("O14").Value = 2.23727414861861
$ws.Range("B15").Value = 0.18842729662083
$ws.Range("D15").Value = 0.2737231428095726
$ws.Range("E15").Value = 0.2014559229930626
$ws.Range("F15").Value = 1.056408717965553
$ws.Range("G15").Value = 0.5064677789801948
$ws.Range("H15").Value = 0.6297682772362521
$ws.Range("J15").Value = 0.1993924502014295
$ws.Range("K15").Value = 0.9608354825908805
$ws.Range("N15").Value = 1.297051230634075
$ws.Range("O15").Value = 2.236915778722505
$ws.Range("B16").Value = 0.1807160757671937
$ws.Range("D16").Value = 0.2685894359366046
$ws.Range("E16").Value = 0.1976952723765919
$ws.Range("F16").Value = 1.052334556669848
$ws.Range("G16").Value = 0.5044624647632787
$ws.Range("H16").Value = 0.6310856628154511
$ws.Range("J16").Value = 0.1956760522741945
$ws.Range("K16").Value = 0.9024507794972862
$ws.Range("N16").Value = 1.302048065881692
$ws.Range("O16").Value = 2.235262201307535
$ws.Range("B17").Value = 0.1760026476314209
$ws.Range("D17").Value = 0.2654663803781148
$ws.Range("E17").Value = 0.195411912137935
$ws.Range("F17").Value = 1.049992282897932
$ws.Range("G17").Value = 0.5033263817142455
$ws.Range("H17").Value = 0.6319702023585307
$ws.Range("J17").Value = 0.1934231921356542
$ws.Range("K17").Value = 0.8665882928403903
$ws.Range("N17").Value = 1.305242895831533
$ws.Range("O17").Value = 2.234604726920963
$ws.Range("B18").Value = 0.1732979801531798
$ws.Range("D18").Value = 0.263679793520339
$ws.Range("E18").Value = 0.1941073063700429
$ws.Range("F18").Value = 1.048702819215492
$ws.Range("G18").Value = 0.5027074984094924
$ws.Range("H18").Value = 0.6325070423794585
$ws.Range("J18").Value = 0.192137349991043
$ws.Range("K18").Value = 0.8459449320058923
$ws.Range("N18").Value = 1.307128084440059
$ws.Range("O18").Value = 2.23435771155539
$ws.Range("B19").Value = 0.1723833269069956
$ws.Range("D19").Value = 0.2630765581622967
$ws.Range("E19").Value = 0.1936670894198826
$ws.Range("F19").Value = 1.048276149453358
$ws.Range("G19").Value = 0.5025038888284996
$ws.Range("H19").Value = 0.6326936304502766
$ws.Range("J19").Value = 0.1917036946801858
$ws.Range("K19").Value = 0.8389527051365349
$ws.Range("N19").Value = 1.307774557331541
$ws.Range("O19").Value = 2.2342965982225
$ws.Range("B20").Value = 0.1765037421617137
$ws.Range("D20").Value = 0.2657978308053828
$ws.Range("E20").Value = 0.1956540774537459
$ws.Range("F20").Value = 1.050235645440139
$ws.Range("G20").Value = 0.5034437421454925
$ws.Range("H20").Value = 0.6318731360948107
$ws.Range("J20").Value = 0.1936619840645193
$ws.Range("K20").Value = 0.8704076035878359
$ws.Range("N20").Value = 1.304897874512513
$ws.Range("O20").Value = 2.234661142372715
$ws.Range("B21").Value = 0.190423431643282
$ws.Range("D21").Value = 0.2750566620915009
$ws.Range("E21").Value = 0.2024341516008548
$ws.Range("F21").Value = 1.057509137099174
$ws.Range("G21").Value = 0.5070146272008458
$ws.Range("H21").Value = 0.6294517538203195
$ws.Range("J21").Value = 0.2003602992693629
$ws.Range("K21").Value = 0.9758948633121918
$ws.Range("N21").Value = 1.29580112889245
$ws.Range("O21").Value = 2.237453075098813
$ws.Range("B22").Value = 0.1995760699531957
$ws.Range("D22").Value = 0.2811933879281128
$ws.Range("E22").Value = 0.2069425406933121
$ws.Range("F22").Value = 1.062776629165512
$ws.Range("G22").Value = 0.5096566884396054
$ws.Range("H22").Value = 0.6281190026070362
$ws.Range("J22").Value = 0.2048263640222814
$ws.Range("K22").Value = 1.044682720650712
$ws.Range("N22").Value = 1.290278959232197
$ws.Range("O22").Value = 2.240446207419097
$ws.Range("B23").Value = 0.1946861421527473
$ws.Range("D23").Value = 0.2779103201489193
$ws.Range("E23").Value = 0.2045292852040603
$ws.Range("F23").Value = 1.059918272957816
$ws.Range("G23").Value = 0.5082183381726537
$ws.Range("H23").Value = 0.62880747602496
$ws.Range("J23").Value = 0.2024346697284614
$ws.Range("K23").Value = 1.00798392891096
$ws.Range("N23").Value = 1.293187574042641
$ws.Range("O23").Value = 2.238741775180927
$ws.Range("B24").Value = 0.1762771812826429
$ws.Range("D24").Value = 0.2656479543367283
$ws.Range("E24").Value = 0.1955445691909432
$ws.Range("F24").Value = 1.050125443241001
$ws.Range("G24").Value = 0.5033905767642466
$ws.Range("H24").Value = 0.6319169315762423
$ws.Range("J24").Value = 0.1935539970854308
$ws.Range("K24").Value = 0.8686809726421245
$ws.Range("N24").Value = 1.305053707759541
$ws.Range("O24").Value = 2.234635228972536
$ws.Range("B25").Value = 0.1566645008407477
$ws.Range("D25").Value = 0.2527928189288247
$ws.Range("E25").Value = 0.186186806094419
$ws.Range("F25").Value = 1.041770351531483
$ws.Range("G25").Value = 0.4995040538771889
$ws.Range("H25").Value = 0.6363445390637708
$ws.Range("J25").Value = 0.1843549906285773
$ws.Range("K25").Value = 0.7178058625325718
$ws.Range("N25").Value = 1.319673914154784
$ws.Range("O25").Value = 2.235219160667413
